$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$gValues = @{
    2 = 1
    3 = 0
    4 = 2
    5 = 1
    6 = 2
    7 = 1
    8 = 1
    9 = 1
    10 = 3
    11 = 2
    12 = 0
    13 = 2
    14 = 2
    15 = 2
    16 = 1
    17 = 2
    18 = 1
    19 = 2
    20 = 0
    21 = 1
    22 = 0
    23 = 3
    24 = 0
    25 = 2
    26 = 1
    27 = 4
    28 = 2
    29 = 2
    30 = 1
    31 = 1
    32 = 0
    33 = 2
    34 = 1
    35 = 1
    36 = 1
    37 = 0
    38 = 2
    39 = 0
    40 = 0
    41 = 0
    42 = 2
    43 = 2
    44 = 1
    45 = 0
    46 = 2
    47 = 1
    48 = 0
    49 = 0
    50 = 1
    51 = 2
    52 = 2
    53 = 1
    55 = 0
    56 = 1
    57 = 1
    58 = 2
    61 = 0
    62 = 0
    63 = 1
    64 = 0
    65 = 3
    66 = 0
    67 = 1
    68 = 2
    69 = 3
}

foreach ($row in $gValues.Keys) {
    $ws.Range("G$row").Value = $gValues[$row]
}
